$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: change A2 text, add B2 as an empty styled cell (reuses existing style 1)
$ws.Range("A2").Value = "Welcome to the Kiot"
$ws.Range("A1").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats

# D3 empty cell with vertical-top alignment (creates style 2)
$ws.Range("D3").VerticalAlignment = -4160    # xlTop

# Row 3: A3 = "Batman" with new font color black (creates style 3: fontId1 + horiz/vert align)
$ws.Range("A3").Value = "Batman"
$ws.Range("A3").Font.Color = 0

# Row 4: A4 = "Welcome to the smartcliff" (moved from A2), same new style
$ws.Range("A4").Value = "Welcome to the smartcliff"
$ws.Range("A4").Font.Color = 0

# Row 5: A5 = "IronMan", same new style
$ws.Range("A5").Value = "IronMan"
$ws.Range("A5").Font.Color = 0

# Column B width (stored width ends up matching Excel's "18" best-fit width)
$ws.Columns("B").ColumnWidth = 17.1666666666667

# Selection
$ws.Range("C4").Select()

# Page setup
$ws.PageSetup.Orientation = 1  # xlPortrait
